$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the three "special" values back to "time"
$ws.Range("F4").Value = "time"
$ws.Range("B9").Value = "time"
$ws.Range("B14").Value = "time"

# Restore the selection (multi-area selection, active cell in C12)
[void]$ws.Range("C7,C12").Select()
[void]$ws.Range("C12").Activate()
